# Refactor CSV comparator output: interlace the file1/file2 columns
# (Name_file1, Name_file2, Age_file1, Age_file2, Location_file1, Location_file2)
# instead of grouping them by file, and drop the now-unused "extra col_file2"
# column that used to sit at the far right.
#
# Current layout : A=ID B=Name_file1 C=Age_file1 D=Location_file1
#                   E=Name_file2 F=Age_file2 G=Location_file2 H=extra col_file2
# Target layout   : A=ID B=Name_file1 C=Name_file2 D=Age_file1 E=Age_file2
#                   F=Location_file1 G=Location_file2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stash the four columns that need to move (values + formatting) in some
# scratch columns far out of the way, so the reshuffle below can't clobber
# data it still needs to read.
$ws.Range("C1:C6").Copy($ws.Range("Z1"))    # Age_file1
$ws.Range("D1:D6").Copy($ws.Range("AA1"))   # Location_file1
$ws.Range("E1:E6").Copy($ws.Range("AB1"))   # Name_file2
$ws.Range("F1:F6").Copy($ws.Range("AC1"))   # Age_file2

# Wipe out the old C:F block (contents + formats) so that pasting shorter /
# blank cells back in doesn't leave stale leftovers behind.
$ws.Range("C1:F6").Clear()

# Paste the stashed columns back in their new, interlaced order.
$ws.Range("AB1:AB6").Copy($ws.Range("C1"))  # Name_file2     -> C
$ws.Range("Z1:Z6").Copy($ws.Range("D1"))    # Age_file1      -> D
$ws.Range("AC1:AC6").Copy($ws.Range("E1"))  # Age_file2      -> E
$ws.Range("AA1:AA6").Copy($ws.Range("F1"))  # Location_file1 -> F

# Get rid of the scratch copies.
$ws.Range("Z1:AC6").Clear()

# Drop the obsolete "extra col_file2" column entirely.
$ws.Columns("H").Clear()
